$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 41 ---
# F41: was a blank-ish text cell, becomes a date value (2020-04-24 -> serial 43945)
# Reuse the existing date format/style (as already used e.g. by F4) rather than
# creating a brand-new number format entry.
$ws.Range("F4").Copy()
$ws.Range("F41").PasteSpecial(-4122)
$ws.Range("F41").Value = "04/24/2020"

# G41: becomes the long "powered by Google..." text (this also updates the shared
# string that used to live in F41, since that text is being replaced here)
$ws.Range("G41").Value = "powered by Google, attention model with pyramid-encoder and decoder structure"

# --- Add a new row to the table (table already includes the blank row 42;
# row 43 is a genuinely new row, so the table range grows by exactly one row) ---
$tbl = $ws.ListObjects.Item(1)
$row43 = $tbl.ListRows.Add()

# Row 42
$ws.Range("A42").Value = "NEURAL MACHINE TRANSLATION BY JOINTLY LEARNING TO ALIGN AND TRANSLATE"
$ws.Range("B42").Value = 2015
$ws.Range("C42").Value = "Dzmitry Bahdanau, KyungHyun Cho, Yoshua Bengio"
$ws.Range("E42").Value = "attention model"
$ws.Range("G42").Value = "Bahdanau Attention model"

# Row 43
$ws.Range("C43").Value = "Minh-Thang Luong, Hieu Pham, Christopher D. Manning"
$ws.Range("A43").Value = "Effective Approaches to Attention-based Neural Machine Translation"
$ws.Range("B43").Value = 2015
$ws.Range("E43").Value = "attention model"

# --- Update selection to match final state ---
$ws.Range("F43").Select()
